$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old column C ("Uang Jalan"), shifting
# C..G -> D..H, to make room for the new "Harga Jual" price column.
$ws.Columns("C").Insert()

# Register the column-level default formatting (right-aligned, like the
# other data columns) for the brand-new column C by stamping a cell well
# outside the used range, then clearing it again -- the engine keeps the
# column-wide default style even after the triggering cell is cleared.
$ws.Range("C10").HorizontalAlignment = -4152
$ws.Range("C10").Clear()

# Column widths (Excel stores width in characters; COM's ColumnWidth is
# ~0.8333 narrower than the stored/display width).
$ws.Columns("C").ColumnWidth = 20.666666666666668
$ws.Columns("D").ColumnWidth = 18.998697916666668

# Copy the header/data-row formatting from the column that used to be C
# (now shifted to D, "Uang Jalan") onto the new C column so the new
# header cell and data cell pick up the same bordered styles.
$ws.Range("D3:D4").Copy()
$ws.Range("C3:C4").PasteSpecial(-4122)

# New header label + value.
$ws.Range("C3").Value = "Harga Jual"

# Selection, matching the author's final cursor position.
$ws.Range("D10").Select()
